# Add missing department to export data
# Replace the "类别" (type) column with a "科室" (department) column

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F previously held the "类别" / ${record.type} data; swap it for
# the department field instead.
$ws.Range("F2").Value = "`${record.department}"
$ws.Range("F1").Value = "科室"

# Reflect the new active selection on the worksheet (was E2, now F2)
$ws.Range("F2").Select()
